# Add "Not assigned" provider rows to the France and Spain provider cost sheets.

$wb = $excel.ActiveWorkbook

# --- France sheet ---
$france = $wb.Worksheets.Item("France")

# Fill in the previously-missing unit cost for Pediatric Medicine (row 10)
$france.Range("C10").Value = 26.5

# Append a new "Not assigned" row
$france.Range("A14").Value = 0
$france.Range("B14").Value = "Not assigned"
$france.Range("C14").Value = 36.839090909090913
$france.Range("C14").NumberFormat = "0.00"
$france.Range("D14").Value = 1

$france.Range("A14").Select()

# --- Spain sheet ---
$spain = $wb.Worksheets.Item("Spain")

# Append a new "Not assigned" row
$spain.Range("A13").Value = 0
$spain.Range("B13").Value = "Not assigned"
$spain.Range("C13").Value = 57.405918307512138
$spain.Range("C13").NumberFormat = "0.00"
$spain.Range("D13").Value = 1

$spain.Activate()
$spain.Range("F21").Select()
